# Add two new rows (26 and 27) of mock-data metadata to the "1D NEW" sheet,
# extending the Table43 table and updating the selection/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1D NEW")

# Row 29
$ws.Range("B29").Value = "1dmockanderrors26.csv"
$ws.Range("C29").Value = 53
$ws.Range("D29").Value = 1000
$ws.Range("E29").Value = 0.3
$ws.Range("F29").Value = 0.1
$ws.Range("G29").Value = 200
$ws.Range("H29").Value = 1
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5
$ws.Range("L29").Value = 1
$ws.Range("M29").Value = "same but with higher central frequency"

# Row 30
$ws.Range("B30").Value = "1dmockanderrors27.csv"
$ws.Range("C30").Value = 53
$ws.Range("D30").Value = 1000
$ws.Range("E30").Value = 0.3
$ws.Range("F30").Value = 0.1
$ws.Range("G30").Value = 200
$ws.Range("H30").Value = 1
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 1
$ws.Range("M30").Value = "noiseless version of 26"

# Grow the table to include the new rows
$table = $ws.ListObjects.Item("Table43")
$table.Resize($ws.Range("B3:M30"))

# Update the selection to match the post-edit state
$ws.Range("M31").Select()
